# Update "想去人数" (interested count, column F) and "最低票价" (min price, column G)
# figures on both the "展览" (sheet1) and "全部类型" (sheet4) worksheets, which hold
# duplicate data. Rows refer to the same records in each sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G2").Value = 50

    $ws.Range("F3").Value = 1068
    $ws.Range("F6").Value = 4586
    $ws.Range("F8").Value = 381
    $ws.Range("F9").Value = 1340
    $ws.Range("F10").Value = 889
    $ws.Range("F12").Value = 963
    $ws.Range("F14").Value = 538
    $ws.Range("F16").Value = 252
}
